$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.481.89"
$ws.Range("E2").Value = "  -1.06%  "

$ws.Range("D3").Value = "'1.910.45"
$ws.Range("E3").Value = "  -1.49%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'239.05"

$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("D7").Value = "'0.4776"
$ws.Range("E7").Value = "  -2.10%  "

$ws.Range("D8").Value = "'0.2832"
$ws.Range("E8").Value = "  -3.50%  "

$ws.Range("D9").Value = "'0.06699"
$ws.Range("E9").Value = "  -2.94%  "

$ws.Range("D10").Value = "'18.63"
$ws.Range("E10").Value = "  -4.88%  "

$ws.Range("D11").Value = "'100.61"
$ws.Range("E11").Value = "  -4.63%  "

$ws.Range("D12").Value = "'1.919.26"
$ws.Range("E12").Value = "  -1.07%  "

$ws.Range("D13").Value = "'0.07678"
$ws.Range("E13").Value = "  -0.94%  "

$ws.Range("D14").Value = "'5.196"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("D15").Value = "'0.6657"
$ws.Range("E15").Value = "  -4.62%  "

$ws.Range("D16").Value = "'30.529.97"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").Value = "'255.33"
$ws.Range("E17").Value = "  -7.32%  "

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("D19").Value = "'0.000007458"
$ws.Range("E19").Value = "  -3.51%  "

$ws.Range("D20").Value = "'12.62"
$ws.Range("E20").Value = "  -3.91%  "

$ws.Range("D21").Value = "'5.367"
$ws.Range("E21").Value = "  -1.30%  "

$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'6.282"
$ws.Range("E23").Value = "  -2.99%  "

$ws.Range("E24").Value = "  -3.93%  "

$ws.Range("D25").Value = "'166.76"
$ws.Range("E25").Value = "  -0.87%  "

$ws.Range("D26").Value = "'19.10"
$ws.Range("E26").Value = "  -2.48%  "

$ws.Range("D27").Value = "'2.049"
$ws.Range("E27").Value = "  -5.21%  "

$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'1.392"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'4.720"
$ws.Range("E29").Value = "  +4.00%  "

$ws.Range("D30").Value = "'0.1004"
$ws.Range("E30").Value = "  -3.50%  "

$ws.Range("D31").Value = "'1.508"
$ws.Range("E31").Value = "  -2.81%  "

$ws.Range("E32").Value = "  -2.79%  "

$ws.Range("D33").Value = "'0.04705"
$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("D34").Value = "'0.7247"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").Value = "'1.104"
$ws.Range("E35").Value = "  -4.31%  "

$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").Value = "'2.705"
$ws.Range("E37").Value = "  -0.76%  "

$ws.Range("D38").Value = "'0.01907"
$ws.Range("E38").Value = "  -4.02%  "

$ws.Range("D39").Value = "'2.605"
$ws.Range("E39").Value = "  -2.16%  "

$ws.Range("D40").Value = "'6.242"
$ws.Range("E40").Value = "  -2.93%  "

$ws.Range("D41").Value = "'74.41"
$ws.Range("E41").Value = "  -3.76%  "

$ws.Range("E42").Value = "  -6.32%  "

$ws.Range("D43").Value = "'0.8596"
$ws.Range("E43").Value = "  -4.47%  "

$ws.Range("D44").Value = "'105.55"
$ws.Range("E44").Value = "  -2.45%  "

$ws.Range("D45").Value = "'1.000"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").Value = "'0.4218"
$ws.Range("E46").Value = "  -4.34%  "

$ws.Range("D47").Value = "'7.334"
$ws.Range("E47").Value = "  -5.07%  "

$ws.Range("D48").Value = "'0.1195"
$ws.Range("E48").Value = "  -3.64%  "

$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'34.63"
$ws.Range("E49").Value = "  -2.94%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'908.17"
$ws.Range("E50").Value = "  -8.37%  "

$ws.Range("D51").Value = "'8.710"
$ws.Range("E51").Value = "  -5.05%  "
